$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D42").Value = "[name=`"W`"]  그녀석이나 그녀석의 비밀이 없더라도 일 진행은 순조롭게 되는 것 같은데 말이지. 용문을 공격하고, 코어를 가동시키고, 그녀석은 전혀 필요 없지 않아?`n"

$ws.Range("D44").Value = "[name=`"탈룰라`"]  코어의 기능은 내가 직접 가동시킬 수 있다. 하지만 리유니온은 정지 방법도 완전히 파악하고 있어야 해. 열쇠의 쓰임새는 바로 거기에 있다.`n"

$ws.Range("D55").Value = "[name=`"W`"]  그래 그래 좋아. 근데 있잖아 리더, 난 열쇠가 대체 몇 개나 있는 건지 알고 싶은데~ 이게 마지막 질문이야, 더는 안 물어볼게!`n"

$ws.Range("D62").Value = "[name=`"W`"]  잠깐, 내가 너무 이것저것 캐물었나? 이렇게 오랜 시간을 들여가며 설명해주니까 왠지 좀 미안한걸.`n"

$ws.Range("D98").Value = "[name=`"W`"]  날 데리고 이렇게나 오래 연극 해주느라 정말 고생이 많았어. 역시 내 눈은 틀리지 않았다니까? 넌 진짜 명배우야.`n"

$ws.Range("D108").Value = "[name=`"W`"]  시작은 니가 했잖아 이 용 뿔 달린 계집애야! 내가 눈에 거슬려서 태워죽이려는 거야? 아니면 네 가면이 벗겨질까 무서워서 입막음 하려는 거야?`n"

$ws.Range("D127").Value = "[name=`"W`"]  그렇게 따지면 니가 더 귀엽지. 으음~ 근데 그 앙증맞은 혀를 잘라내면 훨씬 더 귀여울 것 같지 않아? 혀가 너무 길어 보이지도 않고 말이야.`n"

$ws.Range("D138").Value = "[name=`"W`"]  윽… 일단, 내 취향 열심히 연구해주고 굳이 내 심기를 건드리는 단어까지 골라 써줘서 고맙다고 하면 되는 거지?`n"

$ws.Range("D152").Value = "[name=`"탈룰라`"]  W, W…… 뻔하지 않은가. 내가 네 마음을 알아채려 한 게 아니라, 너 스스로 너의 앙증맞은 비밀을 펼쳐 보인 것이다.`n"

$ws.Range("D156").Value = "타오르는 불길, 사람을 미치게 만드는 타는 냄새, 대량의 조각이 뒤섞인 뜨거운 바람이 내 앞의 모든 사물을 뒤엎는다. 모든 준비는 끝났다. 설령 강인한 살카즈라 해도……`n"

$ws.Range("D174").Value = "숨을 사각지대 하나 없이 사방에서 그녀를 산산조각냈어야 할 폭발이, 모두 속수무책으로 돌아갔다.`n"

$ws.Range("D221").Value = "[name=`"탈룰라`"]  어쩔 수 없었다. 우르수스가 언제 체르노보그를 습격해 우리의 동포를 학살할지 알 수 없었으니, 반드시 그 전에 용문을 함락시켜야 했다.`n"
